# T8_Safety.pptx – slide 3 "folded corner" shape:
# merge the first two paragraphs ("../Sample" / "/Lib/tester")
# into one paragraph ("Lib/tester"), leaving the "[SysB]" paragraph
# untouched.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

try {
    $shape = $s.Shapes.Item("사각형: 모서리가 접힌 도형 8")
} catch {
    # fall back to the known shape index on this slide
    $shape = $s.Shapes.Item(6)
}

$tr = $shape.TextFrame.TextRange

# Current text is "../Sample" + <CR> + "/Lib/tester" + <CR> + "[SysB]".
# Remove "../Sample" together with the paragraph break that follows it,
# as well as the leading "/" of the second paragraph, so the two
# paragraphs collapse into a single "Lib/tester" paragraph.
$lead = $tr.Characters(1, 11)
$lead.Delete()
